$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns on existing header row (L1:N1) ---
$ws.Cells.Item(1, 12).Value = "minimum"
$ws.Cells.Item(1, 13).Value = "maximum"
$ws.Cells.Item(1, 14).Value = "orderAmount"

# --- Row 2 (V00001) gains minimum/maximum/orderAmount ---
$ws.Cells.Item(2, 12).Value = 35
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 40

# --- Row 3: V00002 / Gruppe2 / Vare 2 ---
$ws.Cells.Item(3, 1).Value = "V00002"
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = "Gruppe2"
$ws.Cells.Item(3, 4).Value = "Stk"
$ws.Cells.Item(3, 5).Value = "Vare 2"
$ws.Cells.Item(3, 6).Value = "Varetekst 2"
$ws.Cells.Item(3, 8).Value = 10.6
$ws.Cells.Item(3, 9).Value = 14.4
$ws.Cells.Item(3, 11).Value = "nej"

# --- Row 4: V00003 / Gruppe3 / Vare 3 ---
$ws.Cells.Item(4, 1).Value = "V00003"
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = "Gruppe3"
$ws.Cells.Item(4, 4).Value = "Stk"
$ws.Cells.Item(4, 5).Value = "Vare 3"
$ws.Cells.Item(4, 6).Value = "Varetekst 2"
$ws.Cells.Item(4, 8).Value = 10.6
$ws.Cells.Item(4, 9).Value = 14.4
$ws.Cells.Item(4, 11).Value = "nej"
$ws.Cells.Item(4, 12).Value = 25
$ws.Cells.Item(4, 14).Value = 40

# --- Row 5: V00004 / Diverse / Liter / Vare 4 ---
$ws.Cells.Item(5, 1).Value = "V00004"
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = "Diverse"
$ws.Cells.Item(5, 4).Value = "Liter"
$ws.Cells.Item(5, 5).Value = "Vare 4"
$ws.Cells.Item(5, 6).Value = "Varetekst 2"
$ws.Cells.Item(5, 8).Value = 10.6
$ws.Cells.Item(5, 9).Value = 14.4
$ws.Cells.Item(5, 11).Value = "nej"
$ws.Cells.Item(5, 12).Value = 2500.5
$ws.Cells.Item(5, 12).NumberFormat = "#,##0.00"
$ws.Cells.Item(5, 13).Value = 5000
$ws.Cells.Item(5, 13).NumberFormat = "#,##0"
$ws.Cells.Item(5, 14).Value = 3000

# --- Match the saved selection state ---
$ws.Range("N6").Select() | Out-Null
